# Auto-generated Excel COM-interop script to apply the Ixion_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Cells.Item(13, 8).Value = 5943  # H13 was 5500
$ws.Cells.Item(13, 10).Value = 5943  # J13 was 5500
$ws.Cells.Item(13, 12).Value = 5943  # L13 was 5500
$ws.Cells.Item(13, 14).Value = -6281  # N13 was -5838

# Row 41
$ws.Cells.Item(41, 8).Value = 708.25  # H41 was 765.06665
$ws.Cells.Item(41, 9).Value = 512.7  # I41 was 536.3333
$ws.Cells.Item(41, 10).Value = 903.8  # J41 was 1108.1666
$ws.Cells.Item(41, 11).Value = 512.7  # K41 was 536.3333
$ws.Cells.Item(41, 12).Value = 903.8  # L41 was 1108.1666
$ws.Cells.Item(41, 13).Value = -72.70000000000005  # M41 was -96.33330000000001
$ws.Cells.Item(41, 14).Value = -1783.8  # N41 was -1988.1666

# Row 98
$ws.Cells.Item(98, 8).Value = 1671.5834  # H98 was 1718.2174
$ws.Cells.Item(98, 9).Value = 1767.7894  # I98 was 1832.7222
$ws.Cells.Item(98, 11).Value = 1767.7894  # K98 was 1832.7222
$ws.Cells.Item(98, 13).Value = -269.7893999999999  # M98 was -334.7221999999999

# Row 106
$ws.Cells.Item(106, 8).Value = 95243410  # H106 was 222228560
$ws.Cells.Item(106, 9).Value = 27783146  # I106 was 83340340
$ws.Cells.Item(106, 11).Value = 27783146  # K106 was 83340340
$ws.Cells.Item(106, 13).Value = -27782515  # M106 was -83339709

# Row 122
$ws.Cells.Item(122, 8).Value = 1671.5834  # H122 was 1718.2174
$ws.Cells.Item(122, 9).Value = 1767.7894  # I122 was 1832.7222
$ws.Cells.Item(122, 11).Value = 5303.3682  # K122 was 5498.1666
$ws.Cells.Item(122, 13).Value = -2853.3682  # M122 was -3048.1666

# Row 138
$ws.Cells.Item(138, 8).Value = 4080.011  # H138 was 4162.375
$ws.Cells.Item(138, 9).Value = 1231.9736  # I138 was 1272.8611
$ws.Cells.Item(138, 10).Value = 6122  # J138 was 6162.8076
$ws.Cells.Item(138, 11).Value = 3695.9208  # K138 was 3818.5833
$ws.Cells.Item(138, 12).Value = 18366  # L138 was 18488.4228
$ws.Cells.Item(138, 13).Value = 1444.0792  # M138 was 1321.4167
$ws.Cells.Item(138, 14).Value = -28646  # N138 was -28768.4228

$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Cells.Item(34, 8).Value = 0  # H34 was 12000
$ws.Cells.Item(34, 10).Value = 0  # J34 was 12000
$ws.Cells.Item(34, 12).Value = 0  # L34 was 12000
$ws.Cells.Item(34, 14).ClearContents()  # N34 was -12542

# Row 74
$ws.Cells.Item(74, 8).Value = 1756.5454  # H74 was 1783.3125
$ws.Cells.Item(74, 9).Value = 1573.1111  # I74 was 1612.7059
$ws.Cells.Item(74, 11).Value = 1573.1111  # K74 was 1612.7059
$ws.Cells.Item(74, 13).Value = -699.1111000000001  # M74 was -738.7058999999999

# Row 77
$ws.Cells.Item(77, 8).Value = 1756.5454  # H77 was 1783.3125
$ws.Cells.Item(77, 9).Value = 1573.1111  # I77 was 1612.7059
$ws.Cells.Item(77, 11).Value = 7865.5555  # K77 was 8063.5295
$ws.Cells.Item(77, 13).Value = -3497.5555  # M77 was -3695.5295

# Row 118
$ws.Cells.Item(118, 8).Value = 34800  # H118 was 0
$ws.Cells.Item(118, 10).Value = 34800  # J118 was 0
$ws.Cells.Item(118, 12).Value = 34800  # L118 was 0
$ws.Cells.Item(118, 14).Value = -38114  # N118 was None

$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Cells.Item(40, 8).Value = 50398  # H40 was 50448
$ws.Cells.Item(40, 10).Value = 50398  # J40 was 50448
$ws.Cells.Item(40, 12).Value = 50398  # L40 was 50448
$ws.Cells.Item(40, 14).Value = -50928  # N40 was -50978

# Row 87
$ws.Cells.Item(87, 8).Value = 50000  # H87 was 0
$ws.Cells.Item(87, 10).Value = 50000  # J87 was 0
$ws.Cells.Item(87, 12).Value = 50000  # L87 was 0
$ws.Cells.Item(87, 14).Value = -52496  # N87 was None

# Row 90
$ws.Cells.Item(90, 8).Value = 50000  # H90 was 0
$ws.Cells.Item(90, 10).Value = 50000  # J90 was 0
$ws.Cells.Item(90, 12).Value = 150000  # L90 was 0
$ws.Cells.Item(90, 14).Value = -162480  # N90 was None

# Row 137
$ws.Cells.Item(137, 8).Value = 0  # H137 was 40780
$ws.Cells.Item(137, 10).Value = 0  # J137 was 40780
$ws.Cells.Item(137, 12).Value = 0  # L137 was 40780
$ws.Cells.Item(137, 14).ClearContents()  # N137 was -50980

# Row 140
$ws.Cells.Item(140, 8).Value = 45854.285  # H140 was 42655.555
$ws.Cells.Item(140, 10).Value = 45854.285  # J140 was 42655.555
$ws.Cells.Item(140, 12).Value = 45854.285  # L140 was 42655.555
$ws.Cells.Item(140, 14).Value = -56214.285  # N140 was -53015.555

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Cells.Item(4, 8).Value = 1000  # H4 was 25000
$ws.Cells.Item(4, 9).Value = 1000  # I4 was 25000
$ws.Cells.Item(4, 11).Value = 1000  # K4 was 25000
$ws.Cells.Item(4, 13).Value = -888  # M4 was -24888

# Row 99
$ws.Cells.Item(99, 8).Value = 15648862  # H99 was 12523790
$ws.Cells.Item(99, 10).Value = 25020180  # J99 was 17878272
$ws.Cells.Item(99, 12).Value = 25020180  # L99 was 17878272
$ws.Cells.Item(99, 14).Value = -25023176  # N99 was -17881268

# Row 126
$ws.Cells.Item(126, 8).Value = 15648862  # H126 was 12523790
$ws.Cells.Item(126, 10).Value = 25020180  # J126 was 17878272
$ws.Cells.Item(126, 12).Value = 75060540  # L126 was 53634816
$ws.Cells.Item(126, 14).Value = -75065480  # N126 was -53639756

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Cells.Item(18, 8).Value = 3500  # H18 was 3150
$ws.Cells.Item(18, 9).Value = 3500  # I18 was 3150
$ws.Cells.Item(18, 11).Value = 3500  # K18 was 3150
$ws.Cells.Item(18, 13).Value = -3207  # M18 was -2857

# Row 33
$ws.Cells.Item(33, 8).Value = 0  # H33 was 5800
$ws.Cells.Item(33, 10).Value = 0  # J33 was 5800
$ws.Cells.Item(33, 12).Value = 0  # L33 was 5800
$ws.Cells.Item(33, 14).ClearContents()  # N33 was -6304

# Row 39
$ws.Cells.Item(39, 8).Value = 0  # H39 was 30000
$ws.Cells.Item(39, 10).Value = 0  # J39 was 30000
$ws.Cells.Item(39, 12).Value = 0  # L39 was 30000
$ws.Cells.Item(39, 14).Value = 0  # N39 was -31064

# Row 138
$ws.Cells.Item(138, 8).Value = 45059.43  # H138 was 0
$ws.Cells.Item(138, 10).Value = 45059.43  # J138 was 0
$ws.Cells.Item(138, 12).Value = 45059.43  # L138 was 0
$ws.Cells.Item(138, 14).Value = -55339.43  # N138 was None

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 61521.35  # H7 was 64916.938
$ws.Cells.Item(7, 9).Value = 127301  # I7 was 93052.82000000001
$ws.Cells.Item(7, 10).Value = 3050.5557  # J7 was 3018
$ws.Cells.Item(7, 11).Value = 127301  # K7 was 93052.82000000001
$ws.Cells.Item(7, 12).Value = 3050.5557  # L7 was 3018
$ws.Cells.Item(7, 13).Value = -127189  # M7 was -92940.82000000001
$ws.Cells.Item(7, 14).Value = -3274.5557  # N7 was -3242

# Row 40
$ws.Cells.Item(40, 8).Value = 90913160  # H40 was 90913070
$ws.Cells.Item(40, 9).Value = 200003100  # I40 was 142860320
$ws.Cells.Item(40, 10).Value = 4875.8335  # J40 was 5375
$ws.Cells.Item(40, 11).Value = 200003100  # K40 was 142860320
$ws.Cells.Item(40, 12).Value = 4875.8335  # L40 was 5375
$ws.Cells.Item(40, 13).Value = -200002964  # M40 was -142860184
$ws.Cells.Item(40, 14).Value = -5147.8335  # N40 was -5647

# Row 64
$ws.Cells.Item(64, 8).Value = 32000  # H64 was 30000
$ws.Cells.Item(64, 10).Value = 32000  # J64 was 30000
$ws.Cells.Item(64, 12).Value = 32000  # L64 was 30000
$ws.Cells.Item(64, 14).Value = -32450  # N64 was -30450

# Row 67
$ws.Cells.Item(67, 8).Value = 32000  # H67 was 30000
$ws.Cells.Item(67, 10).Value = 32000  # J67 was 30000
$ws.Cells.Item(67, 12).Value = 32000  # L67 was 30000
$ws.Cells.Item(67, 14).Value = -33560  # N67 was -31560

# Row 68
$ws.Cells.Item(68, 8).Value = 32259858  # H68 was 32259828
$ws.Cells.Item(68, 9).Value = 1771.238  # I68 was 1827
$ws.Cells.Item(68, 10).Value = 100001840  # J68 was 100001630
$ws.Cells.Item(68, 11).Value = 1771.238  # K68 was 1827
$ws.Cells.Item(68, 12).Value = 100001840  # L68 was 100001630
$ws.Cells.Item(68, 13).Value = -1022.238  # M68 was -1078
$ws.Cells.Item(68, 14).Value = -100003338  # N68 was -100003128

# Row 71
$ws.Cells.Item(71, 8).Value = 32259858  # H71 was 32259828
$ws.Cells.Item(71, 9).Value = 1771.238  # I71 was 1827
$ws.Cells.Item(71, 10).Value = 100001840  # J71 was 100001630
$ws.Cells.Item(71, 11).Value = 8856.190000000001  # K71 was 9135
$ws.Cells.Item(71, 12).Value = 500009200  # L71 was 500008150
$ws.Cells.Item(71, 13).Value = -5112.190000000001  # M71 was -5391
$ws.Cells.Item(71, 14).Value = -500016688  # N71 was -500015638

# Row 93
$ws.Cells.Item(93, 8).Value = 125052500  # H93 was 55580104
$ws.Cells.Item(93, 9).Value = 200000  # I93 was 50600
$ws.Cells.Item(93, 10).Value = 166670000  # J93 was 100003710
$ws.Cells.Item(93, 11).Value = 200000  # K93 was 50600
$ws.Cells.Item(93, 12).Value = 166670000  # L93 was 100003710
$ws.Cells.Item(93, 13).Value = -198752  # M93 was -49352
$ws.Cells.Item(93, 14).Value = -166672496  # N93 was -100006206

# Row 94
$ws.Cells.Item(94, 8).Value = 30000  # H94 was 25999.666
$ws.Cells.Item(94, 10).Value = 30000  # J94 was 25999.666
$ws.Cells.Item(94, 12).Value = 30000  # L94 was 25999.666
$ws.Cells.Item(94, 14).Value = -31352  # N94 was -27351.666

# Row 120
$ws.Cells.Item(120, 8).Value = 0  # H120 was 50000
$ws.Cells.Item(120, 10).Value = 0  # J120 was 50000
$ws.Cells.Item(120, 12).Value = 0  # L120 was 50000
$ws.Cells.Item(120, 14).ClearContents()  # N120 was -59676

# Row 126
$ws.Cells.Item(126, 8).Value = 61521.35  # H126 was 64916.938
$ws.Cells.Item(126, 9).Value = 127301  # I126 was 93052.82000000001
$ws.Cells.Item(126, 10).Value = 3050.5557  # J126 was 3018
$ws.Cells.Item(126, 11).Value = 381903  # K126 was 279158.46
$ws.Cells.Item(126, 12).Value = 9151.667099999999  # L126 was 9054
$ws.Cells.Item(126, 13).Value = -379433  # M126 was -276688.46
$ws.Cells.Item(126, 14).Value = -14091.6671  # N126 was -13994

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Cells.Item(63, 8).Value = 33650  # H63 was 40000
$ws.Cells.Item(63, 10).Value = 33650  # J63 was 40000
$ws.Cells.Item(63, 12).Value = 33650  # L63 was 40000
$ws.Cells.Item(63, 14).Value = -34898  # N63 was -41248

# Row 66
$ws.Cells.Item(66, 8).Value = 33650  # H66 was 40000
$ws.Cells.Item(66, 10).Value = 33650  # J66 was 40000
$ws.Cells.Item(66, 12).Value = 100950  # L66 was 120000
$ws.Cells.Item(66, 14).Value = -107190  # N66 was -126240

# Row 126
$ws.Cells.Item(126, 8).Value = 1738.5  # H126 was 1527.2727
$ws.Cells.Item(126, 9).Value = 1301.3334  # I126 was 1132.6666
$ws.Cells.Item(126, 11).Value = 3904.0002  # K126 was 3397.9998
$ws.Cells.Item(126, 13).Value = -1434.0002  # M126 was -927.9998000000001
